$p = $ppt.ActivePresentation

# --- Slide 1: "TextBox 5" title text ---
# "Insights for Creating Successful Films" -> "Insights for Creating Successful Movies"
$s1 = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item(3)
$tr1 = $sh1.TextFrame.TextRange
$run1 = $tr1.Paragraphs(1, 1).Runs(1, 1)
$run1.Text = "Insights for Creating Successful Movies"

# --- Slide 20: "TextBox 4" bullet list ---
$s20 = $p.Slides.Item(20)
$sh20 = $s20.Shapes.Item(10)
$tr20 = $sh20.TextFrame.TextRange

# Paragraph 1: "wise business proposition" -> "wise business venture"
$run20a = $tr20.Paragraphs(1, 1).Runs(1, 1)
$run20a.Text = "In conclusion, this venture stands out to be a wise business venture, as the revenue to be drawn from this is promising in both the domestic and gross numbers. "

# Paragraph 3: remove leading space before "Collaborate"
$run20b = $tr20.Paragraphs(3, 1).Runs(1, 1)
$run20b.Text = "Collaborate with current leading studios for pilots."

# --- Slide 6: "Content Placeholder 3" body text ---
$s6 = $p.Slides.Item(6)
$sh6 = $s6.Shapes.Item(2)
$tr6 = $sh6.TextFrame.TextRange

# Paragraph 5: "As *y  rule of thumb" -> "As my rule of thumb"
$run6a = $tr6.Paragraphs(5, 1).Runs(1, 1)
$run6a.Text = "As my rule of thumb, the columns with missing values that are less than 5%, and the rows with the missing values will be dropped, as this will not significantly impact the analysis."

# Paragraph 7: "The mean was used in case there were outliers..." -> "The median was used, in this case, there were outliers... - to prevent skewness "
$run6b = $tr6.Paragraphs(7, 1).Runs(1, 1)
$run6b.Text = "The median was used, in this case, there were outliers in both datasets - to prevent skewness "
